# feat: add 2022-Q3 data
#
# Duplicate the existing "2022-Q2" sheet (the duplicate is inserted
# immediately before the source sheet), rename the duplicate to
# "2022-Q3", and overwrite its figures with the new Q3 numbers. The
# original "2022-Q2" and "2022-Q1" sheets are left untouched, so they
# keep representing their own quarters. Finally, the "总计" (totals)
# sheet gets a new row for 2022-Q3, with the old rows shifted down.
#
# NOTE: worksheet references returned by this host appear to be bound by
# tab position rather than sheet identity, so after any operation that
# inserts/moves/copies sheets we re-resolve sheets by a fresh
# Worksheets.Item(...) lookup instead of reusing an old variable.

$wb = $excel.ActiveWorkbook

# --- Step 1: duplicate "2022-Q2" (tab index 2) -----------------------
# The copy is placed right before its source, so it lands at index 2
# and pushes the original "2022-Q2" sheet to index 3.
$wb.Worksheets.Item(2).Copy($wb.Worksheets.Item(2))
$wb.Worksheets.Item(2).Name = "2022-Q3"

# --- Step 2: update the new "2022-Q3" sheet's figures ----------------
$q3 = $wb.Worksheets.Item(2)

# Columns D/E/F/G hold text-formatted numbers (e.g. "19.80"); force text
# so trailing zeros survive instead of being normalised away as a number.
$q3.Range("D2:G8").NumberFormat = "@"

# Row 2: 008269 大成睿享混合A
$q3.Range("D2").Value = "19.80"
$q3.Range("E2").Value = "66.91"
$q3.Range("F2").Value = "4.39"
$q3.Range("G2").Value = "0.8692"
$q3.Range("H2").Value = 3

# Row 3: 090013 大成竞争优势混合
$q3.Range("D3").Value = "6.88"
$q3.Range("E3").Value = "61.00"
$q3.Range("F3").Value = "4.16"
$q3.Range("G3").Value = "0.2862"
$q3.Range("H3").Value = 3

# Row 4: 013463 大成致远优势一年持有期混合A
$q3.Range("D4").Value = "3.65"
$q3.Range("E4").Value = "60.88"
$q3.Range("F4").Value = "5.58"
$q3.Range("G4").Value = "0.2037"

# Row 5: 011834, fund renamed 大成投资严选六个月持有期混合型证券投资基金A -> 大成投资严选六月持有混合A
$q3.Range("C5").Value = "大成投资严选六月持有混合A"
$q3.Range("D5").Value = "3.10"
$q3.Range("E5").Value = "66.75"
$q3.Range("F5").Value = "5.80"
$q3.Range("G5").Value = "0.1798"

# Row 6: 008270 大成睿享混合C
$q3.Range("D6").Value = "4.02"
$q3.Range("E6").Value = "66.91"
$q3.Range("F6").Value = "4.39"
$q3.Range("G6").Value = "0.1765"
$q3.Range("H6").Value = 3

# Row 7: 011835, fund renamed 大成投资严选六个月持有期混合型证券投资基金C -> 大成投资严选六月持有混合C
$q3.Range("C7").Value = "大成投资严选六月持有混合C"
$q3.Range("D7").Value = "0.22"
$q3.Range("E7").Value = "66.75"
$q3.Range("F7").Value = "5.80"
$q3.Range("G7").Value = "0.0128"

# Row 8: 013464 大成致远优势一年持有期混合C
$q3.Range("D8").Value = "0.17"
$q3.Range("E8").Value = "60.88"
$q3.Range("F8").Value = "5.58"
$q3.Range("G8").Value = "0.0095"

# --- Step 3: update the "总计" (totals) sheet -------------------------
# Shift old rows down first (row3 -> row4, row2 -> row3), then write the
# new 2022-Q3 row into row 2.
$total = $wb.Worksheets.Item(1)

# Row 4's A cell needs the same style as A2/A3 (bold/border/centered);
# copy it over from A3 before writing the value.
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 7
$total.Range("D4").Value = 1.05

$total.Range("B3").Value = "2022-Q2"
$total.Range("D3").Value = 2.41

$total.Range("B2").Value = "2022-Q3"
$total.Range("D2").Value = 1.74

# --- Step 4: restore the originally-active tab -----------------------
# "2022-Q1" was the selected/active sheet before the edit; keep it that
# way (it has simply moved from tab 3 to tab 4).
$wb.Worksheets.Item("2022-Q1").Activate()
